# Add a second email hyperlink entry below the existing one.
#
# Before: A1 contains "mln02koushik@gmail.com" as a mailto hyperlink.
# After : A2 also contains an email address ("swapnanilsaha26@gmail.com")
#         rendered as plain text (no special hyperlink formatting) but with
#         a live mailto hyperlink attached to it, and the selection moves to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("A2")

# Put the email address text into the cell.
$target.Value = "swapnanilsaha26@gmail.com"

# Attach a mailto hyperlink to the cell, showing the address itself as the
# display text.
$ws.Hyperlinks.Add($target, "mailto:swapnanilsaha26@gmail.com", [Type]::Missing, [Type]::Missing, "swapnanilsaha26@gmail.com")

# Hyperlinks.Add auto-applies Excel's built-in "Hyperlink" style (blue +
# underline). The source file keeps A2 in the default/normal style, so
# restore plain formatting.
$target.Font.Underline = $false
$target.Font.Color = 0

# Match row 2's height to the now-taller text row.
$ws.Rows(2).RowHeight = 14.9

# Move the active selection to the newly filled cell.
$target.Select()
